$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect Price (column D) cells as text before assigning, since some
# new values (e.g. "595.92", "0.999") would otherwise be auto-detected
# by Excel as numbers and lose their original text type.
$dCells = @("D2", "D3", "D5", "D6", "D7", "D12", "D14", "D15", "D16", "D17", "D18", "D21", "D22", "D23", "D25", "D27", "D28", "D30", "D31", "D32", "D34", "D35", "D36", "D37", "D39", "D41", "D42", "D45", "D47", "D48", "D49", "D50", "D51")
foreach ($cellRef in $dCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Column D (Price) updates
$ws.Range("D2").Value = "67.456.32"
$ws.Range("D3").Value = "3.766.92"
$ws.Range("D5").Value = "595.92"
$ws.Range("D6").Value = "169.06"
$ws.Range("D7").Value = "3.765.85"
$ws.Range("D12").Value = "0.457"
$ws.Range("D14").Value = "36.81"
$ws.Range("D15").Value = "4.399.91"
$ws.Range("D16").Value = "3.766.92"
$ws.Range("D17").Value = "19.01"
$ws.Range("D18").Value = "67.547.99"
$ws.Range("D21").Value = "10.54"
$ws.Range("D22").Value = "468.59"
$ws.Range("D23").Value = "0.724"
$ws.Range("D25").Value = "83.90"
$ws.Range("D27").Value = "12.16"
$ws.Range("D28").Value = "10.29"
$ws.Range("D30").Value = "2.91"
$ws.Range("D31").Value = "3.922.45"
$ws.Range("D32").Value = "7.64"
$ws.Range("D34").Value = "30.42"
$ws.Range("D35").Value = "9.18"
$ws.Range("D36").Value = "3.732.49"
$ws.Range("D37").Value = "3.84"
$ws.Range("D39").Value = "5.90"
$ws.Range("D41").Value = "0.999"
$ws.Range("D42").Value = "0.999"
$ws.Range("D45").Value = "8.73"
$ws.Range("D47").Value = "46.25"
$ws.Range("D48").Value = "401.74"
$ws.Range("D49").Value = "0.000275"
$ws.Range("D50").Value = "141.80"
$ws.Range("D51").Value = "0.0354"

# Restore default style so no residual explicit cell style is left behind
foreach ($cellRef in $dCells) {
    $ws.Range($cellRef).Style = "Normal"
}

# Column E (Volume 1h) updates
$ws.Range("E2").Value = "  -1.16%  "
$ws.Range("E3").Value = "  -2.12%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E5").Value = "  -0.87%  "
$ws.Range("E6").Value = "  +0.57%  "
$ws.Range("E7").Value = "  -2.13%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("E10").Value = "  +0.61%  "
$ws.Range("E11").Value = "  +0.49%  "
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("E13").Value = "  +4.73%  "
$ws.Range("E14").Value = "  -0.46%  "
$ws.Range("E15").Value = "  -2.10%  "
$ws.Range("E16").Value = "  -2.11%  "
$ws.Range("E17").Value = "  +5.14%  "
$ws.Range("E18").Value = "  -1.10%  "
$ws.Range("E19").Value = "  -1.39%  "
$ws.Range("E20").Value = "  +0.93%  "
$ws.Range("E21").Value = "  -1.99%  "
$ws.Range("E22").Value = "  +0.55%  "
$ws.Range("E23").Value = "  -1.14%  "
$ws.Range("E24").Value = "  -5.99%  "
$ws.Range("E25").Value = "  +1.24%  "
$ws.Range("E26").Value = "  +0.37%  "
$ws.Range("E27").Value = "  +1.00%  "
$ws.Range("E28").Value = "  +2.86%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("E30").Value = "  -2.05%  "
$ws.Range("E31").Value = "  -1.90%  "
$ws.Range("E32").Value = "  +0.67%  "
$ws.Range("E33").Value = "  -2.46%  "
$ws.Range("E35").Value = "  -3.24%  "
$ws.Range("E36").Value = "  -2.10%  "
$ws.Range("E37").Value = "  +5.52%  "
$ws.Range("E38").Value = "  +0.70%  "
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("E40").Value = "  -2.24%  "
$ws.Range("E41").Value = "  -2.36%  "
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("E43").Value = "  +0.48%  "
$ws.Range("E45").Value = "  +1.54%  "
$ws.Range("E46").Value = "  -0.73%  "
$ws.Range("E47").Value = "  -1.76%  "
$ws.Range("E48").Value = "  -4.36%  "
$ws.Range("E49").Value = "  -6.64%  "
$ws.Range("E50").Value = "  -0.87%  "
$ws.Range("E51").Value = "  -0.83%  "
